$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "1.014") must be forced
# to remain text, matching the workbook convention of storing prices as strings.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "30.715.83"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.124.49"
$ws.Range("E3").Value = "  +1.15%  "
Set-TextValue "D4" "1.014"
$ws.Range("E4").Value = "  +1.06%  "
Set-TextValue "D5" "338.57"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("E6").Value = "  +1.06%  "
Set-TextValue "D7" "0.5270"
$ws.Range("E7").Value = "  +0.88%  "
Set-TextValue "D8" "0.4565"
$ws.Range("E8").Value = "  +1.56%  "
Set-TextValue "D9" "54.95"
$ws.Range("E9").Value = "  +2.10%  "
Set-TextValue "D10" "0.09118"
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("E11").Value = "  +1.87%  "
Set-TextValue "D12" "24.56"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "2.119.66"
$ws.Range("E13").Value = "  +1.17%  "
Set-TextValue "D14" "6.868"
$ws.Range("E14").Value = "  +2.16%  "
Set-TextValue "D15" "8.149"
$ws.Range("E15").Value = "  +5.69%  "
Set-TextValue "D16" "0.00001175"
$ws.Range("E16").Value = "  +4.56%  "
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("E18").Value = "  +1.03%  "
Set-TextValue "D19" "0.06708"
$ws.Range("E19").Value = "  +1.32%  "
Set-TextValue "D20" "19.58"
$ws.Range("E20").Value = "  +2.03%  "
Set-TextValue "D21" "1.012"
$ws.Range("E21").Value = "  +1.09%  "
Set-TextValue "D22" "6.361"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").Value = "30.782.76"
$ws.Range("E23").Value = "  +0.79%  "
Set-TextValue "D24" "13.00"
$ws.Range("E24").Value = "  +5.32%  "
Set-TextValue "D25" "2.364"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").Value = "2.376.37"
$ws.Range("E26").Value = "  +1.57%  "
Set-TextValue "D27" "22.47"
$ws.Range("E27").Value = "  +0.87%  "
Set-TextValue "D28" "165.90"
$ws.Range("E28").Value = "  +1.38%  "
Set-TextValue "D29" "2.564"
$ws.Range("E29").Value = "  -0.53%  "
Set-TextValue "D30" "134.89"
$ws.Range("E30").Value = "  +2.04%  "
Set-TextValue "D31" "1.208"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("E32").Value = "  +0.34%  "
Set-TextValue "D33" "1.664"
$ws.Range("E33").Value = "  -0.28%  "
Set-TextValue "D34" "6.396"
$ws.Range("E34").Value = "  +3.79%  "
Set-TextValue "D35" "3.948"
$ws.Range("E35").Value = "  +1.21%  "
Set-TextValue "D36" "10.63"
$ws.Range("E36").Value = "  +1.38%  "
Set-TextValue "D37" "5.900"
$ws.Range("E37").Value = "  +7.62%  "
Set-TextValue "D38" "0.02664"
$ws.Range("E38").Value = "  +3.65%  "
Set-TextValue "D39" "0.06886"
$ws.Range("E39").Value = "  +1.42%  "
Set-TextValue "D40" "0.2330"
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("E41").Value = "  -0.22%  "
Set-TextValue "D42" "0.6934"
$ws.Range("E42").Value = "  +0.10%  "
Set-TextValue "D43" "1.262"
$ws.Range("E43").Value = "  +0.69%  "
Set-TextValue "D44" "15.19"
$ws.Range("E44").Value = "  +8.59%  "
Set-TextValue "D45" "0.6505"
$ws.Range("E45").Value = "  +2.39%  "
Set-TextValue "D46" "2.316"
$ws.Range("E46").Value = "  +1.88%  "
$ws.Range("E47").Value = "  +18.88%  "
$ws.Range("E48").Value = "  +1.87%  "
Set-TextValue "D50" "83.82"
$ws.Range("E50").Value = "  +2.15%  "
Set-TextValue "D51" "0.07308"
$ws.Range("E51").Value = "  +3.70%  "
